$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Requested quantity" headers on the existing sheets ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B1").Value = "Weekly_PO_Qty"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "PO Forecast"

# Match the outline + page-margin conventions used on the other sheets
$ws3.Outline.SummaryRow = 1
$ws3.Outline.SummaryColumn = 1
$ws3.PageSetup.LeftMargin = 54
$ws3.PageSetup.RightMargin = 54
$ws3.PageSetup.TopMargin = 72
$ws3.PageSetup.BottomMargin = 72
$ws3.PageSetup.HeaderMargin = 36
$ws3.PageSetup.FooterMargin = 36

# --- 3. Header row ---
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# --- 4. Data rows ---
$ws3.Range("A2").Value = 45312.99999999999
$ws3.Range("B2").Value = 34
$ws3.Range("C2").Value = 17.27622103176725
$ws3.Range("D2").Value = 51.06231790689662

$ws3.Range("A3").Value = 45410.99999999999
$ws3.Range("B3").Value = 30
$ws3.Range("C3").Value = 12.91756407548998
$ws3.Range("D3").Value = 46.90032548124486

$ws3.Range("A4").Value = 45417.99999999999
$ws3.Range("B4").Value = 30
$ws3.Range("C4").Value = 12.55393183980752
$ws3.Range("D4").Value = 45.84133881934547

$ws3.Range("A5").Value = 45424.99999999999
$ws3.Range("B5").Value = 29
$ws3.Range("C5").Value = 14.66624327927624
$ws3.Range("D5").Value = 46.36505436219894

$ws3.Range("A6").Value = 45431.99999999999
$ws3.Range("B6").Value = 29
$ws3.Range("C6").Value = 10.5500643617067
$ws3.Range("D6").Value = 44.4485211295066

$ws3.Range("A7").Value = 45501.99999999999
$ws3.Range("B7").Value = 26
$ws3.Range("C7").Value = 7.560579033622414
$ws3.Range("D7").Value = 42.61110947834192

$ws3.Range("A8").Value = 45508.99999999999
$ws3.Range("B8").Value = 26
$ws3.Range("C8").Value = 8.392398542725836
$ws3.Range("D8").Value = 42.3079577025442

$ws3.Range("A9").Value = 45515.99999999999
$ws3.Range("B9").Value = 25
$ws3.Range("C9").Value = 8.586119191881444
$ws3.Range("D9").Value = 40.84281076299107

$ws3.Range("A10").Value = 45522.99999999999
$ws3.Range("B10").Value = 25
$ws3.Range("C10").Value = 9.532898234060976
$ws3.Range("D10").Value = 41.15694485093424

$ws3.Range("A11").Value = 45529.99999999999
$ws3.Range("B11").Value = 25
$ws3.Range("C11").Value = 9.193504971496758
$ws3.Range("D11").Value = 41.11029261256859

$ws3.Range("A12").Value = 45536.99999999999
$ws3.Range("B12").Value = 24
$ws3.Range("C12").Value = 7.484592513195003
$ws3.Range("D12").Value = 40.31679235733174

$ws3.Range("A13").Value = 45543.99999999999
$ws3.Range("B13").Value = 24
$ws3.Range("C13").Value = 7.248104100728908
$ws3.Range("D13").Value = 40.64080696871279

$ws3.Range("A14").Value = 45550.99999999999
$ws3.Range("B14").Value = 24
$ws3.Range("C14").Value = 6.540674083602029
$ws3.Range("D14").Value = 39.19623149789996

$ws3.Range("A15").Value = 45557.99999999999
$ws3.Range("B15").Value = 24
$ws3.Range("C15").Value = 7.003083350403943
$ws3.Range("D15").Value = 40.6086065981058

$ws3.Range("A16").Value = 45564.99999999999
$ws3.Range("B16").Value = 23
$ws3.Range("C16").Value = 5.938116744767477
$ws3.Range("D16").Value = 39.42531152955241

$ws3.Range("A17").Value = 45571.99999999999
$ws3.Range("B17").Value = 23
$ws3.Range("C17").Value = 7.613071750379506
$ws3.Range("D17").Value = 39.80182110168974

# --- 5. Formatting: match the style used for header row / date column on the other sheets ---
$ws1.Range("A1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

$ws1.Range("A2").Copy()
$ws3.Range("A2:A17").PasteSpecial(-4122)  # xlPasteFormats
